$p = $ppt.ActivePresentation
$tm = [char]0x2122

# ---------------------------------------------------------------------------
# 1. New slide inserted at position 2: "Data Analysis Techniques"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Data Analysis Techniques`t"

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "Used the "
$r = $body2.InsertAfter("DataAnalyzinator")
$r = $r.InsertAfter(" 3000$tm created by McNichols, Ramos, and Enid")
$r = $r.InsertAfter("`nCalculated GPS location of each suspect using phone bearing and ")
$r = $r.InsertAfter("lat")
$r = $r.InsertAfter("/long changes over time")
$r = $r.InsertAfter("`nGrabbed any users that had locations overlapping by 250m or less at the same hour")

# A slide is transiently created and removed here in between -- matches the
# slide-id bookkeeping of the source deck (ids 256,260,257,259,258,262,263:
# id 261 is consumed and released before the two trailing slides are added).
$dummy = $p.Slides.Add($p.Slides.Count + 1, 2)

# ---------------------------------------------------------------------------
# 2. New slide appended: "Calculating Correlation of Activity + Location"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Add($p.Slides.Count + 1, 2)

$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Calculating Correlation of Activity + Location`t"

$body6 = $s6.Shapes.Item(2).TextFrame.TextRange
$body6.Text = "Used 3 separate correlation coefficient calculations from "
$r = $body6.InsertAfter("scipy")
$r = $r.InsertAfter(":")
$r = $r.InsertAfter("`nSpearman")
$r = $r.InsertAfter("`nKendall")
$r = $r.InsertAfter("`nPearson")
$r = $r.InsertAfter("`nCorrelation between latitude and activity less than .18")
$r = $r.InsertAfter("`nCorrelation between longitude and activity less than .07")
$r = $r.InsertAfter("`nConclusion: No correlation between activity and latitude or longitude separately")

$body6.Paragraphs(2).IndentLevel = 2
$body6.Paragraphs(3).IndentLevel = 2
$body6.Paragraphs(4).IndentLevel = 2

$dummy.Delete()

# ---------------------------------------------------------------------------
# 3. New slide appended: "Calculating Correlation of Activity + Timestamp"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Add($p.Slides.Count + 1, 2)

$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Calculating Correlation of Activity + Timestamp`t"

$body7 = $s7.Shapes.Item(2).TextFrame.TextRange
$body7.Text = "Used same 3 correlation coefficient calculations from "
$r = $body7.InsertAfter("scipy")
$r = $r.InsertAfter(":")
$r = $r.InsertAfter("`nSpearman")
$r = $r.InsertAfter("`nKendall")
$r = $r.InsertAfter("`nPearson")
$r = $r.InsertAfter("`nCorrelation between time and activity .96, .69, .86")
$r = $r.InsertAfter("`nConclusion: Positive correlation between activity and time")

$body7.Paragraphs(2).IndentLevel = 2
$body7.Paragraphs(3).IndentLevel = 2
$body7.Paragraphs(4).IndentLevel = 2
